$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.634.69"
$ws.Range("E2").Value = "  +2.21%  "
$ws.Range("D3").Value = "1.870.33"
$ws.Range("E3").Value = "  +2.18%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4627"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3885"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07874"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9752"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "1.870.62"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.997"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.703"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06972"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.12"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001003"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.83"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("D21").Value = "28.646.94"
$ws.Range("E21").Value = "  +2.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.284"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.01"
$ws.Range("D23").ClearFormats()
$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").Value = "2.169.17"
$ws.Range("E24").Value = "  +5.83%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.117"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.78"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.19"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.798"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.40"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09352"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9140"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.265"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.330"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05784"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02108"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.54%  "
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.779"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5634"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1787"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.773"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.73"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5313"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.156"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.142"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.830"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "113.23"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.402"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.56%  "
$ws.Range("E51").Value = "  +0.31%  "
